$wb = $excel.ActiveWorkbook
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # New row 5: the event that was previously in row 4 (青田 event) shifts down one row.
    $ws.Range("A5").Value = 4
    $ws.Range("B5").NumberFormat = "@"
    $ws.Range("B5").Value = "2024-10-02"
    $ws.Range("B5").Style = "Normal"
    $ws.Range("C5").Value = "青田·未闻展名国漫嘉年华"
    $ws.Range("D5").Value = "瓯南街道百悦城4幢 西娜君澜大饭店"
    $ws.Range("E5").Value = "2024.10.02 09:00-10.02 17:00"
    $ws.Range("F5").Value = 58
    $ws.Range("G5").Value = 45
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=91328"
    $ws.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202408/w8uKtdlg1724147282076.jpeg"

    $ws.Range("A3").Copy()
    $ws.Range("A5").PasteSpecial(-4122)

    # Row 4: new event (丽水·熙梦动漫游戏展).
    $ws.Range("A4").Value = 3
    $ws.Range("B4").NumberFormat = "@"
    $ws.Range("B4").Value = "2024-10-01"
    $ws.Range("B4").Style = "Normal"
    $ws.Range("C4").Value = "丽水·熙梦动漫游戏展"
    $ws.Range("D4").Value = "城北街798号 莱茵体育生活馆"
    $ws.Range("E4").Value = "2024.10.01 10:00-10.01 17:00"
    $ws.Range("F4").Value = 0
    $ws.Range("G4").Value = 45
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=92235"
    $ws.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202409/JHGyuq6R1725938726636.jpeg"

    $ws.Range("A3").Copy()
    $ws.Range("A6").PasteSpecial(-4122)

    # Row 6: new event (缙云·星辰动漫游戏展嘉年华).
    $ws.Range("A6").Value = 5
    $ws.Range("B6").NumberFormat = "@"
    $ws.Range("B6").Value = "2024-10-03"
    $ws.Range("B6").Style = "Normal"
    $ws.Range("C6").Value = "缙云·星辰动漫游戏展嘉年华"
    $ws.Range("D6").Value = "黄龙路38号 中意大酒店(缙云店)"
    $ws.Range("E6").Value = "2024.10.03 10:00-10.03 17:00"
    $ws.Range("F6").Value = 0
    $ws.Range("G6").Value = 45
    $ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=92236"
    $ws.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202409/S13hVYA01725280725848.jpeg"
}
